$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 334
$ws.Range("F2").Value = 45200
$ws.Range("G2").Value = 30042
$ws.Range("H2").Value = 45231
$ws.Range("E3").Value = 29983
$ws.Range("F3").Value = 45170
$ws.Range("E4").Value = 29983
$ws.Range("F4").Value = 45170
$ws.Range("G4").Value = 30042
$ws.Range("H4").Value = 45231
$ws.Range("E5").Value = 30011
$ws.Range("F5").Value = 45200
$ws.Range("G5").Value = 30042
$ws.Range("H5").Value = 45231
$ws.Range("C6").Value = 442
$ws.Range("F6").Value = 45170
$ws.Range("G6").Value = 30042
$ws.Range("H6").Value = 45231
$ws.Range("E7").Value = 29952
$ws.Range("F7").Value = 45139
$ws.Range("G7").Value = 30042
$ws.Range("H7").Value = 45231
$ws.Range("D8").Value = 409
$ws.Range("E8").Value = 29983
$ws.Range("F8").Value = 45170
$ws.Range("H8").Value = 45231
$ws.Range("E9").Value = 29983
$ws.Range("F9").Value = 45170
$ws.Range("G9").Value = 30042
$ws.Range("H9").Value = 45231
$ws.Range("C10").Value = 500
$ws.Range("D10").Value = 482
$ws.Range("E10").Value = 29983
$ws.Range("F10").Value = 45170
$ws.Range("H10").Value = 45231
$ws.Range("E11").Value = 29952
$ws.Range("F11").Value = 45139
$ws.Range("G11").Value = 30042
$ws.Range("H11").Value = 45231
$ws.Range("C12").Value = 370
$ws.Range("D12").Value = 351
$ws.Range("F12").Value = 45170
$ws.Range("H12").Value = 45231
$ws.Range("C13").Value = 466
$ws.Range("F13").Value = 45170
$ws.Range("G13").Value = 30042
$ws.Range("H13").Value = 45231
$ws.Range("C14").Value = 423
$ws.Range("D14").Value = 395
$ws.Range("F14").Value = 45170
$ws.Range("H14").Value = 45231
$ws.Range("C15").Value = 383
$ws.Range("F15").Value = 45170
$ws.Range("G15").Value = 30011
$ws.Range("H15").Value = 45231
$ws.Range("C16").Value = 454
$ws.Range("D16").Value = 409
$ws.Range("F16").Value = 45170
$ws.Range("H16").Value = 45231
$ws.Range("C17").Value = 369
$ws.Range("D17").Value = 393
$ws.Range("F17").Value = 45170
$ws.Range("H17").Value = 45231
$ws.Range("E18").Value = 29983
$ws.Range("F18").Value = 45170
$ws.Range("G18").Value = 30042
$ws.Range("H18").Value = 45231
$ws.Range("D19").Value = 397
$ws.Range("E19").Value = 29983
$ws.Range("F19").Value = 45170
$ws.Range("H19").Value = 45231
$ws.Range("C20").Value = 480
$ws.Range("F20").Value = 45170
$ws.Range("G20").Value = 30042
$ws.Range("H20").Value = 45231
$ws.Range("C21").Value = 309
$ws.Range("F21").Value = 45170
$ws.Range("G21").Value = 30042
$ws.Range("H21").Value = 45231
$ws.Range("C22").Value = 322
$ws.Range("D22").Value = 366
$ws.Range("F22").Value = 45170
$ws.Range("H22").Value = 45231
$ws.Range("C23").Value = 259
$ws.Range("D23").Value = 401
$ws.Range("F23").Value = 45139
$ws.Range("H23").Value = 45231
$ws.Range("C24").Value = 215
$ws.Range("D24").Value = 409
$ws.Range("F24").Value = 45200
$ws.Range("H24").Value = 45231
$ws.Range("D25").Value = 314
$ws.Range("E25").Value = 29952
$ws.Range("F25").Value = 45139
$ws.Range("H25").Value = 45231
$ws.Range("C26").Value = 332
$ws.Range("D26").Value = 312
$ws.Range("F26").Value = 45170
$ws.Range("H26").Value = 45231
$ws.Range("E27").Value = 29983
$ws.Range("F27").Value = 45170
$ws.Range("G27").Value = 30042
$ws.Range("H27").Value = 45231
$ws.Range("D28").Value = 379
$ws.Range("E28").Value = 29952
$ws.Range("F28").Value = 45139
$ws.Range("H28").Value = 45231
$ws.Range("C29").Value = 261
$ws.Range("D29").Value = 234
$ws.Range("F29").Value = 45170
$ws.Range("H29").Value = 45231
$ws.Range("D30").Value = 216
$ws.Range("E30").Value = 29983
$ws.Range("F30").Value = 45170
$ws.Range("H30").Value = 45231
$ws.Range("C31").Value = 393
$ws.Range("F31").Value = 45170
$ws.Range("G31").Value = 30042
$ws.Range("H31").Value = 45231
$ws.Range("E32").Value = 29983
$ws.Range("F32").Value = 45170
$ws.Range("G32").Value = 30042
$ws.Range("H32").Value = 45231
$ws.Range("C33").Value = 454
$ws.Range("D33").Value = 397
$ws.Range("F33").Value = 45200
$ws.Range("H33").Value = 45231
$ws.Range("C34").Value = 201
$ws.Range("D34").Value = 316
$ws.Range("F34").Value = 45170
$ws.Range("H34").Value = 45231
$ws.Range("C35").Value = 403
$ws.Range("D35").Value = 316
$ws.Range("F35").Value = 45139
$ws.Range("H35").Value = 45231
$ws.Range("C36").Value = 500
$ws.Range("D36").Value = 409
$ws.Range("E36").Value = 30011
$ws.Range("F36").Value = 45200
$ws.Range("H36").Value = 45231
$ws.Range("C37").Value = 465
$ws.Range("D37").Value = 316
$ws.Range("F37").Value = 45170
$ws.Range("H37").Value = 45231
$ws.Range("C38").Value = 358
$ws.Range("D38").Value = 366
$ws.Range("F38").Value = 45170
$ws.Range("H38").Value = 45231
$ws.Range("C39").Value = 226
$ws.Range("D39").Value = 222
$ws.Range("F39").Value = 45170
$ws.Range("H39").Value = 45231
$ws.Range("C40").Value = 285
$ws.Range("D40").Value = 313
$ws.Range("F40").Value = 45170
$ws.Range("H40").Value = 45231
$ws.Range("C41").Value = 393
$ws.Range("D41").Value = 315
$ws.Range("F41").Value = 45170
$ws.Range("H41").Value = 45231
$ws.Range("C42").Value = 238
$ws.Range("D42").Value = 217
$ws.Range("F42").Value = 45170
$ws.Range("H42").Value = 45231
$ws.Range("C43").Value = 472
$ws.Range("D43").Value = 316
$ws.Range("F43").Value = 45170
$ws.Range("H43").Value = 45231
$ws.Range("C44").Value = 404
$ws.Range("D44").Value = 303
$ws.Range("F44").Value = 45139
$ws.Range("H44").Value = 45231
$ws.Range("C45").Value = 375
$ws.Range("D45").Value = 316
$ws.Range("F45").Value = 45139
$ws.Range("H45").Value = 45231
$ws.Range("C46").Value = 333
$ws.Range("D46").Value = 297
$ws.Range("F46").Value = 45170
$ws.Range("H46").Value = 45231
$ws.Range("C47").Value = 334
$ws.Range("D47").Value = 259
$ws.Range("F47").Value = 45170
$ws.Range("H47").Value = 45231
$ws.Range("C48").Value = 357
$ws.Range("D48").Value = 314
$ws.Range("F48").Value = 45170
$ws.Range("H48").Value = 45231
$ws.Range("C49").Value = 296
$ws.Range("D49").Value = 312
$ws.Range("F49").Value = 45139
$ws.Range("H49").Value = 45231
$ws.Range("C50").Value = 358
$ws.Range("D50").Value = 237
$ws.Range("F50").Value = 45170
$ws.Range("H50").Value = 45231
$ws.Range("D51").Value = 316
$ws.Range("E51").Value = 29799
$ws.Range("F51").Value = 45139
$ws.Range("H51").Value = 45231
$ws.Range("C52").Value = 345
$ws.Range("D52").Value = 314
$ws.Range("F52").Value = 45170
$ws.Range("H52").Value = 45231
